# Apply crypto price/volume updates from the Mon Apr 15 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.271.14"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "3.102.84"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.26%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.095.22"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000218"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.63%  "
$ws.Range("D15").Value = "3.611.35"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "63.332.94"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "3.113.17"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "506.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.07%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -10.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.55%  "
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "526.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.02%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0414"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("D39").Value = "3.080.77"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0794"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("E42").Value = "  -10.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +78.06%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.253"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.65%  "
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Value = "0.0₃0509"
$ws.Range("E51").Value = "  -8.28%  "
